$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-22 correspond to Generation (column B) 0..20 -> Fitness (column C) becomes 7320
$ws.Range("C2:C22").Value = 7320

# Rows 23-52 correspond to Generation (column B) 21..50 -> Fitness (column C) becomes 7295
$ws.Range("C23:C52").Value = 7295
